# Add data for 2022-08-28
# Updates the 2022 year-to-date totals (column I) in the "Citywide Totals" and
# "By Neighborhood" summary sheets, plus the corresponding per-neighborhood
# detail sheets, to reflect one additional day of violent-crime records.
# A handful of 2019 (column F) figures are also corrected where the source
# data was retroactively reclassified.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 4736
$ws.Range("I3").Value = 4924
$ws.Range("F4").Value = 1864
$ws.Range("I4").Value = 1129
$ws.Range("I5").Value = 453
$ws.Range("I6").Value = 5353
$ws.Range("F7").Value = 24053
$ws.Range("I7").Value = 16595

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 130
$ws.Range("I4").Value = 64
$ws.Range("I5").Value = 53
$ws.Range("I7").Value = 531
$ws.Range("I8").Value = 1016
$ws.Range("I10").Value = 117
$ws.Range("I11").Value = 253
$ws.Range("I15").Value = 188
$ws.Range("I19").Value = 465
$ws.Range("I20").Value = 404
$ws.Range("I22").Value = 45
$ws.Range("I23").Value = 159
$ws.Range("I29").Value = 1057
$ws.Range("I31").Value = 160
$ws.Range("I33").Value = 764
$ws.Range("I34").Value = 80
$ws.Range("I35").Value = 21
$ws.Range("I36").Value = 224
$ws.Range("I37").Value = 527
$ws.Range("I42").Value = 568
$ws.Range("I43").Value = 134
$ws.Range("I47").Value = 115
$ws.Range("I54").Value = 368
$ws.Range("I55").Value = 182
$ws.Range("I59").Value = 29
$ws.Range("I60").Value = 87
$ws.Range("F63").Value = 156
$ws.Range("I64").Value = 149
$ws.Range("F65").Value = 456
$ws.Range("I65").Value = 369
$ws.Range("I67").Value = 652
$ws.Range("I68").Value = 61
$ws.Range("I71").Value = 45
$ws.Range("I73").Value = 142
$ws.Range("I74").Value = 32
$ws.Range("I77").Value = 101
$ws.Range("I78").Value = 233
$ws.Range("I79").Value = 461
$ws.Range("I83").Value = 346
$ws.Range("I85").Value = 752
$ws.Range("I88").Value = 152
$ws.Range("I90").Value = 205
$ws.Range("I91").Value = 193
$ws.Range("I94").Value = 158
$ws.Range("I95").Value = 271
$ws.Range("I99").Value = 310
$ws.Range("F101").Value = 24053
$ws.Range("I101").Value = 16595

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 199
$ws.Range("I3").Value = 303
$ws.Range("I7").Value = 752

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I4").Value = 22
$ws.Range("I6").Value = 65
$ws.Range("I7").Value = 253

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 321
$ws.Range("I3").Value = 283
$ws.Range("I5").Value = 27
$ws.Range("I6").Value = 324
$ws.Range("I7").Value = 1016

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I2").Value = 178
$ws.Range("I3").Value = 165
$ws.Range("I7").Value = 531

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 169
$ws.Range("I3").Value = 170
$ws.Range("I7").Value = 527

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I2").Value = 83
$ws.Range("I7").Value = 310

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I3").Value = 230
$ws.Range("I6").Value = 211
$ws.Range("I7").Value = 652

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I2").Value = 47
$ws.Range("I3").Value = 44
$ws.Range("I7").Value = 160

$ws = $wb.Worksheets.Item('New City')
$ws.Range("F4").Value = 34
$ws.Range("I6").Value = 109
$ws.Range("F7").Value = 456
$ws.Range("I7").Value = 369

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I3").Value = 131
$ws.Range("I7").Value = 346

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I6").Value = 50
$ws.Range("I7").Value = 271

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I3").Value = 282
$ws.Range("I4").Value = 34
$ws.Range("I6").Value = 240
$ws.Range("I7").Value = 764

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I6").Value = 177
$ws.Range("I7").Value = 368

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 312
$ws.Range("I3").Value = 367
$ws.Range("I6").Value = 285
$ws.Range("I7").Value = 1057

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 172
$ws.Range("I6").Value = 128
$ws.Range("I7").Value = 465

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 149
$ws.Range("I4").Value = 47
$ws.Range("I7").Value = 568

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("I3").Value = 22
$ws.Range("I6").Value = 52
$ws.Range("I7").Value = 117

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I4").Value = 30
$ws.Range("I7").Value = 233

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I2").Value = 54
$ws.Range("I6").Value = 57
$ws.Range("I7").Value = 182

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I3").Value = 58
$ws.Range("I6").Value = 44
$ws.Range("I7").Value = 159

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I2").Value = 61
$ws.Range("I6").Value = 55
$ws.Range("I7").Value = 193

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 135
$ws.Range("I7").Value = 461

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("I2").Value = 43
$ws.Range("I4").Value = 9
$ws.Range("I7").Value = 149

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I3").Value = 122
$ws.Range("I4").Value = 30
$ws.Range("I7").Value = 404

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I3").Value = 72
$ws.Range("I7").Value = 224

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("I2").Value = 34
$ws.Range("I7").Value = 80

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I6").Value = 89
$ws.Range("I7").Value = 158

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I3").Value = 37
$ws.Range("I7").Value = 115

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I3").Value = 44
$ws.Range("I4").Value = 14
$ws.Range("I7").Value = 188

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("I6").Value = 6
$ws.Range("I7").Value = 21

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("I3").Value = 48
$ws.Range("I7").Value = 142

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("I6").Value = 8
$ws.Range("I7").Value = 29

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I4").Value = 14
$ws.Range("I7").Value = 130

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I2").Value = 45
$ws.Range("I6").Value = 43
$ws.Range("I7").Value = 152

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("I2").Value = 13
$ws.Range("I7").Value = 53

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I3").Value = 46
$ws.Range("I7").Value = 205

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("I2").Value = 22
$ws.Range("I7").Value = 61

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("I5").Value = 6
$ws.Range("I7").Value = 87

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I4").Value = 10
$ws.Range("I7").Value = 134

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("I6").Value = 14
$ws.Range("I7").Value = 45

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("I2").Value = 15
$ws.Range("I7").Value = 45

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("I3").Value = 36
$ws.Range("I7").Value = 101

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 64

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("I3").Value = 2
$ws.Range("I7").Value = 32
